# Applies updated ifo GDP component analysis preprocessing values
# to the error-table summary rows (Q2-Q9 blocks, rows 2-10).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = @(0.08696016040226752, 0.2244626349384274, 0.06217355075105529, 0.2493462467153963, 0.24189316101004, 15)
    3  = @(0.3031818451638543, 0.3412803276417196, 0.1783459041844743, 0.4223101990059846, 0.3050817571715199, 14)
    4  = @(0.4726833070249565, 0.4978299971988559, 0.4530240298790372, 0.6730705979903128, 0.4987257742493065, 13)
    5  = @(0.6055399787264154, 0.6086332724655883, 0.5906225627394526, 0.7685197217635034, 0.4942695954492832, 12)
    6  = @(0.5882471822302764, 0.5882471822302764, 0.4329084706557275, 0.6579578030966177, 0.3091295773284307, 11)
    7  = @(0.4671934630864089, 0.4696202125542849, 0.2831907857527659, 0.5321567304401644, 0.2685786741804888, 10)
    8  = @(0.3964777679756588, 0.4142371153650239, 0.2111306574106052, 0.4594895618081059, 0.2463291325150163, 9)
    9  = @(0.430140519527954, 0.430140519527954, 0.2113118878782925, 0.4596867279771002, 0.1776210167919809, 6)
    10 = @(0.3316630433627284, 0.3316630433627284, 0.1188298149235608, 0.344717007012362, 0.1150832780485532, 3)
}

foreach ($row in $values.Keys) {
    $rowVals = $values[$row]
    $ws.Cells.Item($row, 2).Value = $rowVals[0]  # B
    $ws.Cells.Item($row, 3).Value = $rowVals[1]  # C
    $ws.Cells.Item($row, 4).Value = $rowVals[2]  # D
    $ws.Cells.Item($row, 5).Value = $rowVals[3]  # E
    $ws.Cells.Item($row, 6).Value = $rowVals[4]  # F
    $ws.Cells.Item($row, 7).Value = $rowVals[5]  # G
}
